# Update "Elapsed Duration(Hrs)" values (column G) across the R1, R2, R4, R5, R6
# sheets to reflect the report being re-generated ~12m28s (748s) later.
$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Sheet = "R1"; Cell = "G2"; Value = "3923:39:59" },
    @{ Sheet = "R1"; Cell = "G3"; Value = "63:12:37" },
    @{ Sheet = "R2"; Cell = "G2"; Value = "12105:03:39" },
    @{ Sheet = "R2"; Cell = "G3"; Value = "3234:47:08" },
    @{ Sheet = "R2"; Cell = "G4"; Value = "472:58:42" },
    @{ Sheet = "R4"; Cell = "G2"; Value = "2950:53:28" },
    @{ Sheet = "R4"; Cell = "G3"; Value = "178:05:43" },
    @{ Sheet = "R5"; Cell = "G2"; Value = "424:52:27" },
    @{ Sheet = "R6"; Cell = "G2"; Value = "65:24:45" }
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    $ws.Range($u.Cell).Value = $u.Value
}
